$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend header row 1 with two new columns (P, Q), continuing the 0..15 sequence
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15

# Copy the header style (bold, bordered, centered) from O1 onto the new header cells
$ws.Range("O1").Copy()
$ws.Range("P1:Q1").PasteSpecial(-4122)

# Update data columns I, K, M and O for rows 2-25 (values flip between 1 and 2)
$ws.Range("I2:I25").Value = 2
$ws.Range("K2:K25").Value = 1
$ws.Range("M2:M25").Value = 2
$ws.Range("O2:O25").Value = 1

# Add new data columns P and Q for rows 2-25
$ws.Range("P2:P25").Value = 2
$ws.Range("Q2:Q25").Value = 2
